$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Individuals")
$ws.Range("F6").Value = 117
